# Apply the crypto price/volume table refresh described in the commit
# "Updated cryptos list on Sat Apr 13 22:36:58 UTC 2024 with GitHub Actions".
#
# The sheet stores everything in columns B:E as plain text. Price (D) and
# Volume (E) values often look like numbers (e.g. "0.996", "6.01"), so they
# are written with a leading apostrophe - exactly like typing them into Excel
# by hand - to force a text cell instead of letting Excel reinterpret them as
# numeric values. Coin names (B) and links (C) never look numeric, so they are
# assigned directly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'63.573.72"
$ws.Range("E2").Value = "'  -5.08%  "

# Row 3
$ws.Range("D3").Value = "'2.977.98"
$ws.Range("E3").Value = "'  -7.50%  "

# Row 4
$ws.Range("D4").Value = "'0.996"
$ws.Range("E4").Value = "'  -0.45%  "

# Row 5
$ws.Range("D5").Value = "'550.98"
$ws.Range("E5").Value = "'  -7.80%  "

# Row 6
$ws.Range("D6").Value = "'136.54"
$ws.Range("E6").Value = "'  -9.49%  "

# Row 7
$ws.Range("E7").Value = "'  -0.39%  "

# Row 8
$ws.Range("D8").Value = "'2.951.11"
$ws.Range("E8").Value = "'  -8.64%  "

# Row 9
$ws.Range("D9").Value = "'0.472"
$ws.Range("E9").Value = "'  -13.30%  "

# Row 10
$ws.Range("D10").Value = "'0.149"
$ws.Range("E10").Value = "'  -14.56%  "

# Row 11
$ws.Range("D11").Value = "'6.01"
$ws.Range("E11").Value = "'  -8.50%  "

# Row 12
$ws.Range("D12").Value = "'0.442"
$ws.Range("E12").Value = "'  -11.17%  "

# Row 13
$ws.Range("D13").Value = "'33.26"
$ws.Range("E13").Value = "'  -14.55%  "

# Row 14
$ws.Range("D14").Value = "'0.0000209"
$ws.Range("E14").Value = "'  -15.03%  "

# Row 15
$ws.Range("D15").Value = "'3.398.77"
$ws.Range("E15").Value = "'  -9.16%  "

# Row 16
$ws.Range("D16").Value = "'62.992.95"
$ws.Range("E16").Value = "'  -6.06%  "

# Row 17
$ws.Range("E17").Value = "'  -4.91%  "

# Row 18
$ws.Range("D18").Value = "'2.937.72"
$ws.Range("E18").Value = "'  -8.83%  "

# Row 19
$ws.Range("D19").Value = "'475.91"
$ws.Range("E19").Value = "'  -10.44%  "

# Row 20
$ws.Range("D20").Value = "'6.30"
$ws.Range("E20").Value = "'  -12.12%  "

# Row 21
$ws.Range("D21").Value = "'13.05"
$ws.Range("E21").Value = "'  -12.19%  "

# Row 22
$ws.Range("D22").Value = "'0.646"
$ws.Range("E22").Value = "'  -15.01%  "

# Row 23
$ws.Range("D23").Value = "'6.69"
$ws.Range("E23").Value = "'  -15.56%  "

# Row 24
$ws.Range("B24").Value = "InternetComputer(DFINITY)"
$ws.Range("C24").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D24").Value = "'12.19"
$ws.Range("E24").Value = "'  -11.65%  "

# Row 25
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").Value = "'76.40"
$ws.Range("E25").Value = "'  -10.51%  "

# Row 27
$ws.Range("D27").Value = "'2.63"
$ws.Range("E27").Value = "'  -17.35%  "

# Row 28
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").Value = "'7.30"
$ws.Range("E28").Value = "'  -10.12%  "

# Row 29
$ws.Range("B29").Value = "ImmutableX"
$ws.Range("C29").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D29").Value = "'1.94"
$ws.Range("E29").Value = "'  -12.12%  "

# Row 30
$ws.Range("D30").Value = "'25.09"
$ws.Range("E30").Value = "'  -13.98%  "

# Row 31
$ws.Range("B31").Value = "Stacks"
$ws.Range("C31").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D31").Value = "'2.45"
$ws.Range("E31").Value = "'  -7.24%  "

# Row 32
$ws.Range("B32").Value = "Mantle"
$ws.Range("C32").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D32").Value = "'1.09"
$ws.Range("E32").Value = "'  -4.43%  "

# Row 33
$ws.Range("D33").Value = "'0.991"
$ws.Range("E33").Value = "'  -1.13%  "

# Row 34
$ws.Range("D34").Value = "'495.46"
$ws.Range("E34").Value = "'  -8.70%  "

# Row 35
$ws.Range("D35").Value = "'50.95"
$ws.Range("E35").Value = "'  -4.84%  "

# Row 36
$ws.Range("D36").Value = "'5.65"
$ws.Range("E36").Value = "'  -13.11%  "

# Row 37
$ws.Range("D37").Value = "'4.92"
$ws.Range("E37").Value = "'  -13.59%  "

# Row 38
$ws.Range("D38").Value = "'0.0395"
$ws.Range("E38").Value = "'  -7.34%  "

# Row 39
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "'0.0764"
$ws.Range("E39").Value = "'  -11.47%  "

# Row 40
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").Value = "'0.116"
$ws.Range("E40").Value = "'  -7.98%  "

# Row 41
$ws.Range("D41").Value = "'8.00"
$ws.Range("E41").Value = "'  -14.46%  "

# Row 42
$ws.Range("D42").Value = "'2.749.66"
$ws.Range("E42").Value = "'  -5.82%  "

# Row 44
$ws.Range("D44").Value = "'2.35"
$ws.Range("E44").Value = "'  -10.50%  "

# Row 45
$ws.Range("D45").Value = "'0.228"
$ws.Range("E45").Value = "'  -13.32%  "

# Row 46
$ws.Range("B46").Value = "PEPE"
$ws.Range("C46").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D46").Value = "'0.0₃0501"
$ws.Range("E46").Value = "'  -14.20%  "

# Row 47
$ws.Range("B47").Value = "Fetch.AI"
$ws.Range("C47").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D47").Value = "'1.91"
$ws.Range("E47").Value = "'  -9.71%  "

# Row 48
$ws.Range("B48").Value = "Monero"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D48").Value = "'114.22"
$ws.Range("E48").Value = "'  -6.81%  "

# Row 49
$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").Value = "'0.102"
$ws.Range("E49").Value = "'  -9.82%  "

# Row 50
$ws.Range("D50").Value = "'22.34"
$ws.Range("E50").Value = "'  -15.98%  "

# Row 51
$ws.Range("D51").Value = "'1.94"
$ws.Range("E51").Value = "'  -19.25%  "
